$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "data": add column U (header "30. 11. 2021"), fill U2:U349, fix a
# handful of previously-mis-keyed T column values, and bump the footer date.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# New header cell U1, formatted like the other header cells (copy format from T1)
$ws1.Cells.Item(1, 21).Value = "30. 11. 2021"
$ws1.Range("T1").Copy()
$ws1.Range("U1").PasteSpecial(-4122)

# Fix a few T-column values that were corrected in this revision
$tFixesSheet1 = @{252=0.04;253=0.21;254=0.71;259=0.135;264=0.09;265=0.095}
foreach ($row in $tFixesSheet1.Keys) {
    $ws1.Cells.Item($row, 20).Value = $tFixesSheet1[$row]
}

# New column U data values (rows 2-349)
$uValuesSheet1 = @{2=0.77;3=0.005;4=0.02;5=0.06;6=0.05;7=0.095;8=0.68;9=0.015;10=0.03;11=0.11;12=0.07000000000000001;13=0.095;14=0.73;15=0.005;16=0.02;17=0.06;18=0.07000000000000001;19=0.115;20=0.86;21=0;22=0.02;23=0.04;24=0.02;25=0.06;26=0.72;27=0.01;28=0.03;29=0.09;30=0.05;31=0.1;32=0.76;33=0.01;34=0.02;35=0.06;36=0.07000000000000001;37=0.08;38=0.9;39=0.005;40=0.015;41=0.02;42=0.03;43=0.03;44=0.63;45=0.015;46=0.03;47=0.12;48=0.07000000000000001;49=0.135;50=0.6899999999999999;51=0.01;52=0.02;53=0.07000000000000001;54=0.09;55=0.12;56=0.88;57=0.005;58=0.015;59=0.02;60=0.04;61=0.04;62=0.82;63=0;64=0.03;65=0.05;66=0.02;67=0.08;68=0.91;69=0;70=0.005;71=0.02;72=0.01;73=0.055;74=0.92;75=0.005;76=0.015;77=0.02;78=0.02;79=0.02;80=0.77;81=0.01;82=0.02;83=0.06;84=0.05;85=0.09;86=0.96;87=0;88=0.005;89=0.015;90=0.015;91=0.005;92=0.68;93=0;94=0.005;95=0.1;96=0.09;97=0.125;98=0.58;99=0.005;100=0.005;101=0.1;102=0.13;103=0.18;104=0.77;105=0.005;106=0.02;107=0.07000000000000001;108=0.05;109=0.08500000000000001;110=0.7;111=0;112=0.03;113=0.05;114=0.12;115=0.1;116=0.77;117=0;118=0.01;119=0.07000000000000001;120=0.09;121=0.06;122=0.78;123=0.01;124=0.03;125=0.06;126=0.04;127=0.08;128=0.8100000000000001;129=0.005;130=0.02;131=0.04;132=0.05;133=0.075;134=0.74;135=0.01;136=0.02;137=0.08;138=0.05;139=0.1;140=0.72;141=0.01;142=0.03;143=0.09;144=0.07000000000000001;145=0.08;146=0.78;147=0.005;148=0.02;149=0.04;150=0.06;151=0.095;152=0.9;153=0.005;154=0.015;155=0.01;156=0.02;157=0.05;158=0.64;159=0.03;160=0.015;161=0.12;162=0.08;163=0.115;164=0.6899999999999999;165=0.005;166=0.03;167=0.09;168=0.08;169=0.105;170=0.83;171=0;172=0.03;173=0.06;174=0.02;175=0.06;176=0.76;177=0.01;178=0.02;179=0.07000000000000001;180=0.04;181=0.1;182=0.76;183=0.005;184=0.03;185=0.07000000000000001;186=0.05;187=0.08500000000000001;188=0.74;189=0.005;190=0.015;191=0.05;192=0.1;193=0.09;194=0.85;195=0;196=0.015;197=0.03;198=0.05;199=0.055;200=0.7;201=0.015;202=0.02;203=0.1;204=0.06;205=0.105;206=0.71;207=0.005;208=0.04;209=0.08;210=0.06;211=0.105;212=0.7;213=0;214=0.015;215=0.05;216=0.13;217=0.105;218=0.8;219=0;220=0.02;221=0.05;222=0.07000000000000001;223=0.06;224=0.87;225=0;226=0.02;227=0.04;228=0.01;229=0.06;230=0.83;231=0;232=0.03;233=0.04;234=0.03;235=0.07000000000000001;236=0.8100000000000001;237=0.01;238=0.01;239=0.07000000000000001;240=0.04;241=0.06;242=0.9;243=0.005;244=0.01;245=0.01;246=0.015;247=0.06;248=0.61;249=0.03;250=0.015;251=0.11;252=0.07000000000000001;253=0.165;254=0.76;255=0.005;256=0.04;257=0.04;258=0.05;259=0.105;260=0.82;261=0.005;262=0.015;263=0.06;264=0.05;265=0.05;266=0.83;267=0;268=0.01;269=0.05;270=0.05;271=0.06;272=0.72;273=0.015;274=0.04;275=0.08;276=0.04;277=0.105;278=0.86;279=0;280=0.03;281=0.04;282=0.05;283=0.02;284=0.76;285=0.005;286=0.015;287=0.05;288=0.07000000000000001;289=0.1;290=0.89;291=0.005;292=0.02;293=0.04;294=0.02;295=0.025;296=0.75;297=0;298=0.015;299=0.05;300=0.07000000000000001;301=0.115;302=0.73;303=0.005;304=0.01;305=0.07000000000000001;306=0.09;307=0.095;308=0.8100000000000001;309=0.005;310=0.03;311=0.04;312=0.04;313=0.075;314=0.76;315=0.01;316=0.02;317=0.08;318=0.04;319=0.09;320=0.7;321=0.02;322=0.03;323=0.08;324=0.09;325=0.08;326=0.79;327=0.005;328=0.02;329=0.06;330=0.04;331=0.08500000000000001;332=0.68;333=0;334=0.02;335=0.05;336=0.09;337=0.16;338=0.74;339=0.01;340=0.03;341=0.1;342=0.05;343=0.07000000000000001;344=0.83;345=0.01;346=0.02;347=0.05;348=0.04;349=0.05}
foreach ($row in $uValuesSheet1.Keys) {
    $ws1.Cells.Item($row, 21).Value = $uValuesSheet1[$row]
}

# Footer / title row text bump
$ws1.Cells.Item(350, 1).Value = "Život během pandemie, Zájem o bezplatné očkování proti koronaviru, % respondentů celkově a ve skupinách, aktualizace 8. 12. 2021"

# ---------------------------------------------------------------------------
# Sheet "pocetR": add column T (header "30. 11. 2021"), fill T2:T59, fix a
# handful of previously-mis-keyed S column values, and bump the footer date.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# New header cell T1, formatted like the other header cells (copy format from S1)
$ws2.Cells.Item(1, 20).Value = "30. 11. 2021"
$ws2.Range("S1").Copy()
$ws2.Range("T1").PasteSpecial(-4122)

# Fix a few S-column values that were corrected in this revision
$sFixesSheet2 = @{43=184;44=621;45=537;46=252}
foreach ($row in $sFixesSheet2.Keys) {
    $ws2.Cells.Item($row, 19).Value = $sFixesSheet2[$row]
}

# New column T data values (rows 2-59)
$tValuesSheet2 = @{2=1790;3=452;4=647;5=691;6=829;7=608;8=353;9=430;10=426;11=243;12=399;13=182;14=110;15=777;16=59;17=97;18=85;19=829;20=144;21=113;22=704;23=872;24=918;25=242;26=319;27=311;28=210;29=328;30=380;31=929;32=419;33=206;34=236;35=580;36=259;37=131;38=130;39=349;40=160;41=75;42=106;43=204;44=633;45=585;46=266;47=331;48=256;49=274;50=268;51=424;52=517;53=610;54=663;55=322;56=1468;57=339;58=653;59=798}
foreach ($row in $tValuesSheet2.Keys) {
    $ws2.Cells.Item($row, 20).Value = $tValuesSheet2[$row]
}

# Footer / title row text bump
$ws2.Cells.Item(60, 1).Value = "Život během pandemie, Zájem o bezplatné očkování proti koronaviru, velikost dotázaného souboru celkově a ve skupinách, aktualizace 8. 12. 2021"

# Trailing empty placeholder cell T60 (mirrors the blank inlineStr cells B60:S60)
$ws2.Cells.Item(60, 20).NumberFormat = "General"

Write-Output "edit complete"
